# Regenerate save_data: recompute the "K" column (column G) values.
# (Originally this column tracked a "Strike#" style stat; it has been
#  regenerated using the new K-based calculation, and the updated
#  s_vals are written back into the sheet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 2
    6  = 1
    7  = 3
    8  = 1
    9  = 2
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 0
    15 = 2
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 2
    25 = 1
    26 = 1
    27 = 1
    30 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
